# Apply the "new .ttl from Google sheet has been generated" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: ConceptScheme URI
$ws.Range("B1").Value = "http://ontology.deic.dk/cv/DTUbib"

# Row 3: PREFIX URI
$ws.Range("C3").Value = "http://ontology.deic.dk/cv/DTUbib"

# Row 10: dct:title
$ws.Range("B10").Value = "Vocab DTUbib"

# Row 11: dct:description
$ws.Range("B11").Value = "Test Vocabulary for M4M workshop"

# Row 19: vars:Test -> vars:test, clear description
$ws.Range("A19").Value = "vars:test"
$ws.Range("B19").Value = "test"
$ws.Range("E19").Value = ""

# Row 20: vars:Computerscientist -> vars:test2, clear description and broader
$ws.Range("A20").Value = "vars:test2"
$ws.Range("B20").Value = "test2"
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""

# Row 21: vars:Computerscience -> vars:, clear label and description
$ws.Range("A21").Value = "vars:"
$ws.Range("B21").Value = ""
$ws.Range("E21").Value = ""
